$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85, shifting existing rows 85:194 down to 86:195
$ws.Rows("85:85").Insert()

# Populate the new row 85 with the latest price record
$ws.Range("A85").Value = 8
$ws.Range("B85").Value = "Terminal La Palmera de La Serena"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 45225
$ws.Range("D85").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = 100114007
$ws.Range("G85").Value = "Jengibre"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 400
$ws.Range("K85").Value = 28000
$ws.Range("L85").Value = 29000
$ws.Range("M85").Value = 28500
$ws.Range("N85").Value = "`$/caja 13 kilos"
$ws.Range("O85").Value = "Perú"
$ws.Range("P85").Value = 2192
$ws.Range("Q85").Value = 13
$ws.Range("R85").Value = "Hortaliza"
